$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Dia"
$ws.Range("B3").Value = "Dio"
$ws.Range("B4").Value = "Dpo"

$ws.Range("D2").Value = "joee@yopmail.com"
$ws.Range("D3").Value = "adile@yopmail.com"
$ws.Range("D4").Value = "nabil@yopmail.com"

$ws.Range("E13").Select()
